# Initial CRUD work on Horses
# Adds new schema rows (Color/Markings/Height/About/BestFriend and
# Active/ActiveComments) to the CCR Horse Inventory "Horse" schema block,
# and highlights the existing core string fields (Name/RegisteredName/
# BirthYear/Breed/Gender) with a light-green fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new rows right before the existing "Picture" row (row 9)
# to document the new Color / Markings / Height / About / BestFriend
# fields on the Horse object. ---
$ws.Rows.Item(9).Resize(5).Insert()

$ws.Range("A9").Value = "String"
$ws.Range("B9").Value = "Color"

$ws.Range("A10").Value = "String"
$ws.Range("B10").Value = "Markings"

$ws.Range("A11").Value = "Number"
$ws.Range("B11").Value = "Height"

$ws.Range("A12").Value = "String"
$ws.Range("B12").Value = "About"

$ws.Range("A13").Value = "String"
$ws.Range("B13").Value = "BestFriend"

# --- Insert 2 new rows right after the "Picture" row (now row 14) to
# document the new Active / ActiveComments fields. ---
$ws.Rows.Item(15).Resize(2).Insert()

$ws.Range("A15").Value = "Boolean"
$ws.Range("B15").Value = "Active"

$ws.Range("A16").Value = "String"
$ws.Range("B16").Value = "ActiveComments"

# --- Highlight the core identity fields (Name, RegisteredName,
# BirthYear, Breed, Gender) with a light-green fill. ---
$ws.Range("A4:C8").Interior.Color = 5296274

# --- Leave the selection on the highlighted block, matching the
# author's last on-screen selection. ---
$ws.Range("A4:C8").Select()
